$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1191032.2
$ws.Range("J17").Value = 1191032.2
$ws.Range("L17").Value = 3573096.6
$ws.Range("N17").Value = -3573432.6
$ws.Range("H43").Value = 967.94116
$ws.Range("I43").Value = 624.5714
$ws.Range("J43").Value = 1208.3
$ws.Range("K43").Value = 624.5714
$ws.Range("L43").Value = 1208.3
$ws.Range("M43").Value = -555.5714
$ws.Range("N43").Value = -1346.3
$ws.Range("H81").Value = 27900
$ws.Range("J81").Value = 27900
$ws.Range("L81").Value = 27900
$ws.Range("N81").Value = -29896
$ws.Range("H84").Value = 27900
$ws.Range("J84").Value = 27900
$ws.Range("L84").Value = 83700
$ws.Range("N84").Value = -93684
$ws.Range("H96").Value = 19232680
$ws.Range("I96").Value = 41668476
$ws.Range("J96").Value = 1998.1428
$ws.Range("K96").Value = 125005428
$ws.Range("L96").Value = 5994.428400000001
$ws.Range("M96").Value = -125004055
$ws.Range("N96").Value = -8740.428400000001
$ws.Range("H100").Value = 10754550
$ws.Range("I100").Value = 14493867
$ws.Range("J100").Value = 4013.875
$ws.Range("K100").Value = 14493867
$ws.Range("L100").Value = 4013.875
$ws.Range("M100").Value = -14493326
$ws.Range("N100").Value = -5095.875
$ws.Range("H112").Value = 1170.5454
$ws.Range("J112").Value = 1170.9333
$ws.Range("L112").Value = 3512.7999
$ws.Range("N112").Value = -5728.7999
$ws.Range("H115").Value = 761.5
$ws.Range("H131").Value = 4533.778
$ws.Range("I131").Value = 279.85715
$ws.Range("J131").Value = 4857.446
$ws.Range("K131").Value = 839.5714499999999
$ws.Range("L131").Value = 14572.338
$ws.Range("M131").Value = 4200.428550000001
$ws.Range("N131").Value = -24652.338
$ws.Range("H138").Value = 2606.693
$ws.Range("I138").Value = 1271.9434
$ws.Range("J138").Value = 4627.8857
$ws.Range("K138").Value = 3815.8302
$ws.Range("L138").Value = 13883.6571
$ws.Range("M138").Value = 1324.1698
$ws.Range("N138").Value = -24163.6571

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3059.84
$ws.Range("I32").Value = 2711.1458
$ws.Range("J32").Value = 11428.5
$ws.Range("K32").Value = 2711.1458
$ws.Range("L32").Value = 11428.5
$ws.Range("M32").Value = -2424.1458
$ws.Range("N32").Value = -12002.5
$ws.Range("H45").Value = 1505.2778
$ws.Range("I45").Value = 1191.4546
$ws.Range("K45").Value = 1191.4546
$ws.Range("M45").Value = -814.4546
$ws.Range("H52").Value = 41166.668
$ws.Range("J52").Value = 41166.668
$ws.Range("L52").Value = 41166.668
$ws.Range("N52").Value = -41802.668
$ws.Range("H61").Value = 687.4400000000001
$ws.Range("I61").Value = 632.75
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 632.75
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -420.75
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 3388.4783
$ws.Range("I74").Value = 3465.9546
$ws.Range("J74").Value = 1684
$ws.Range("K74").Value = 3465.9546
$ws.Range("L74").Value = 1684
$ws.Range("M74").Value = -2591.9546
$ws.Range("N74").Value = -3432
$ws.Range("H77").Value = 3388.4783
$ws.Range("I77").Value = 3465.9546
$ws.Range("J77").Value = 1684
$ws.Range("K77").Value = 17329.773
$ws.Range("L77").Value = 8420
$ws.Range("M77").Value = -12961.773
$ws.Range("N77").Value = -17156
$ws.Range("H104").Value = 31112.5
$ws.Range("J104").Value = 31112.5
$ws.Range("L104").Value = 31112.5
$ws.Range("N104").Value = -38100.5
$ws.Range("H136").Value = 687.4400000000001
$ws.Range("I136").Value = 632.75
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 1898.25
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 651.75
$ws.Range("N136").Value = -11100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2532.7407
$ws.Range("I80").Value = 671.9286
$ws.Range("J80").Value = 4536.6924
$ws.Range("K80").Value = 671.9286
$ws.Range("L80").Value = 4536.6924
$ws.Range("M80").Value = 326.0714
$ws.Range("N80").Value = -6532.6924
$ws.Range("H83").Value = 2532.7407
$ws.Range("I83").Value = 671.9286
$ws.Range("J83").Value = 4536.6924
$ws.Range("K83").Value = 3359.643
$ws.Range("L83").Value = 22683.462
$ws.Range("M83").Value = 1632.357
$ws.Range("N83").Value = -32667.462
$ws.Range("H126").Value = 35966.668
$ws.Range("J126").Value = 35966.668
$ws.Range("L126").Value = 35966.668
$ws.Range("N126").Value = -45846.668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1712.9474
$ws.Range("I132").Value = 721.21875
$ws.Range("K132").Value = 2163.65625
$ws.Range("M132").Value = 366.34375
$ws.Range("H138").Value = 56865
$ws.Range("J138").Value = 56865
$ws.Range("L138").Value = 56865
$ws.Range("N138").Value = -67145
$ws.Range("H139").Value = 41033.332
$ws.Range("J139").Value = 41033.332
$ws.Range("L139").Value = 41033.332
$ws.Range("N139").Value = -51313.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 274.47058
$ws.Range("I33").Value = 158
$ws.Range("J33").Value = 323
$ws.Range("K33").Value = 948
$ws.Range("L33").Value = 1938
$ws.Range("M33").Value = -665
$ws.Range("N33").Value = -2504
$ws.Range("H138").Value = 1617.2084
$ws.Range("I138").Value = 1177.2727
$ws.Range("J138").Value = 1989.4615
$ws.Range("K138").Value = 3531.8181
$ws.Range("L138").Value = 5968.3845
$ws.Range("M138").Value = 1608.1819
$ws.Range("N138").Value = -16248.3845
$ws.Range("H139").Value = 1615.9259
$ws.Range("J139").Value = 2064.2856
$ws.Range("L139").Value = 6192.8568
$ws.Range("N139").Value = -16472.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1548.2222
$ws.Range("I132").Value = 1275.4054
$ws.Range("J132").Value = 2810
$ws.Range("K132").Value = 3826.2162
$ws.Range("L132").Value = 8430
$ws.Range("M132").Value = -1296.2162
$ws.Range("N132").Value = -13490
$ws.Range("H138").Value = 40150
$ws.Range("J138").Value = 40150
$ws.Range("L138").Value = 40150
$ws.Range("N138").Value = -50430
$ws.Range("H139").Value = 45444.5
$ws.Range("J139").Value = 45444.5
$ws.Range("L139").Value = 45444.5
$ws.Range("N139").Value = -55724.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3547.1333
$ws.Range("I136").Value = 4108.5835
$ws.Range("J136").Value = 1301.3334
$ws.Range("K136").Value = 12325.7505
$ws.Range("L136").Value = 3904.0002
$ws.Range("M136").Value = -9775.750499999998
$ws.Range("N136").Value = -9004.0002
$ws.Range("H138").Value = 47963.547
$ws.Range("J138").Value = 47963.547
$ws.Range("L138").Value = 47963.547
$ws.Range("N138").Value = -58243.547

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 29066
$ws.Range("J94").Value = 29066
$ws.Range("L94").Value = 29066
$ws.Range("N94").Value = -30868
$ws.Range("H131").Value = 29440
$ws.Range("J131").Value = 29440
$ws.Range("L131").Value = 29440
$ws.Range("N131").Value = -39520
